$d = $word.ActiveDocument

# 1) "June" -> "Sept" in the "October 2015 - June 2018" line (first occurrence only)
$rng = $d.Content
$found = $rng.Find.Execute("June 2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $juneRange = $d.Range($rng.Start, $rng.Start + 4)
    $juneRange.Text = "Sept"
}

# 2) Move the "_GoBack" bookmark from its old spot (just before "Amherst") to
#    right after the newly-edited " 2018" text.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Sept 2018", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $gobackRange = $d.Range($rng2.End, $rng2.End)
    $d.Bookmarks.Add("_GoBack", $gobackRange)
}
